$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (date) values between the two sets of rows.
# Rows 2 and 3 move from 2021-12-29 (44559) to 2022-01-13 (44574)
$ws.Range("D2").Value = 44574
$ws.Range("D3").Value = 44574

# Rows 6 and 7 move from 2022-01-13 (44574) to 2021-12-29 (44559)
$ws.Range("D6").Value = 44559
$ws.Range("D7").Value = 44559
